{"js": "// Replace the arithmetic expression in every cell of the first (only)\n// table in the document body, in row-major reading order, with the new\n// set of expressions from the target revision.\nconst newValues = [\n  [\"74-69=\", \"2+0=\", \"2+86=\", \"14+36=\", \"88-49=\"],\n  [\"1+23=\", \"17+79=\", \"40-8=\", \"65-20=\", \"22+3=\"],\n  [\"95-47=\", \"20+62=\", \"98-28=\", \"61+13=\", \"7+59=\"],\n  [\"26-4=\", \"25+4=\", \"93-62=\", \"7+87=\", \"80-15=\"],\n  [\"16+72=\", \"13+48=\", \"98-77=\", \"57-8=\", \"54-52=\"],\n  [\"43+31=\", \"60+1=\", \"1+36=\", \"9+38=\", \"57-5=\"],\n  [\"24+5=\", \"54+21=\", \"30-12=\", \"97-41=\", \"88-76=\"],\n  [\"52+38=\", \"4+14=\", \"81+9=\", \"62-40=\", \"70-19=\"],\n  [\"14-12=\", \"36+1=\", \"10+7=\", \"97-16=\", \"71+27=\"],\n  [\"99-49=\", \"21+53=\", \"12+38=\", \"93-42=\", \"34-24=\"],\n  [\"56+28=\", \"97-73=\", \"39-9=\", \"79-43=\", \"63-9=\"],\n  [\"76-59=\", \"17+8=\", \"33-16=\", \"7+46=\", \"48+49=\"],\n  [\"77-1=\", \"24+49=\", \"73-9=\", \"60-10=\", \"58-47=\"],\n  [\"33+4=\", \"82-39=\", \"27+71=\", \"69+30=\", \"8+29=\"],\n  [\"77-28=\", \"96+1=\", \"38-37=\", \"55-32=\", \"96-36=\"],\n  [\"5+62=\", \"1+17=\", \"17+73=\", \"15+76=\", \"94-90=\"],\n  [\"19-8=\", \"31+30=\", \"9+67=\", \"13-12=\", \"4+32=\"],\n  [\"85-27=\", \"26+62=\", \"43-5=\", \"61+22=\", \"97-76=\"],\n  [\"89+10=\", \"45+51=\", \"55+14=\", \"47+44=\", \"86-27=\"],\n  [\"97-75=\", \"65+19=\", \"2+41=\", \"37+52=\", \"18+1=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[r][c];\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the arithmetic expression in every cell of the first (only)\n# table in the document body, in row-major reading order, with the new\n# set of expressions from the target revision.\n$newValues = @(\n    @(\"74-69=\", \"2+0=\", \"2+86=\", \"14+36=\", \"88-49=\"),\n    @(\"1+23=\", \"17+79=\", \"40-8=\", \"65-20=\", \"22+3=\"),\n    @(\"95-47=\", \"20+62=\", \"98-28=\", \"61+13=\", \"7+59=\"),\n    @(\"26-4=\", \"25+4=\", \"93-62=\", \"7+87=\", \"80-15=\"),\n    @(\"16+72=\", \"13+48=\", \"98-77=\", \"57-8=\", \"54-52=\"),\n    @(\"43+31=\", \"60+1=\", \"1+36=\", \"9+38=\", \"57-5=\"),\n    @(\"24+5=\", \"54+21=\", \"30-12=\", \"97-41=\", \"88-76=\"),\n    @(\"52+38=\", \"4+14=\", \"81+9=\", \"62-40=\", \"70-19=\"),\n    @(\"14-12=\", \"36+1=\", \"10+7=\", \"97-16=\", \"71+27=\"),\n    @(\"99-49=\", \"21+53=\", \"12+38=\", \"93-42=\", \"34-24=\"),\n    @(\"56+28=\", \"97-73=\", \"39-9=\", \"79-43=\", \"63-9=\"),\n    @(\"76-59=\", \"17+8=\", \"33-16=\", \"7+46=\", \"48+49=\"),\n    @(\"77-1=\", \"24+49=\", \"73-9=\", \"60-10=\", \"58-47=\"),\n    @(\"33+4=\", \"82-39=\", \"27+71=\", \"69+30=\", \"8+29=\"),\n    @(\"77-28=\", \"96+1=\", \"38-37=\", \"55-32=\", \"96-36=\"),\n    @(\"5+62=\", \"1+17=\", \"17+73=\", \"15+76=\", \"94-90=\"),\n    @(\"19-8=\", \"31+30=\", \"9+67=\", \"13-12=\", \"4+32=\"),\n    @(\"85-27=\", \"26+62=\", \"43-5=\", \"61+22=\", \"97-76=\"),\n    @(\"89+10=\", \"45+51=\", \"55+14=\", \"47+44=\", \"86-27=\"),\n    @(\"97-75=\", \"65+19=\", \"2+41=\", \"37+52=\", \"18+1=\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Count; $c++) {\n        $cell = $tbl.Cell($r + 1, $c + 1)\n        $cell.Range.Text = $row[$c]\n    }\n}\n"}
